$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 4: Admin / Admin321 / Invalid credentials
$ws.Range("A4").Value = "Admin "
$ws.Range("B4").Value = "Admin321"
$ws.Range("C4").Value = "Invalid credentials"

# New row 5: abhi / 1234 / Invalid credentials
$ws.Range("A5").Value = "abhi"
$ws.Range("B5").Value = 1234
$ws.Range("C5").Value = "Invalid credentials"

# Reuse the existing "Invalid credentials" cell style (wrap text, Consolas 10pt)
# from C2 instead of letting Excel synthesize a brand-new style.
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C4:C5").PasteSpecial(-4122) | Out-Null

# Match the row heights used by the other data rows.
$ws.Rows.Item(4).RowHeight = 26.25
$ws.Rows.Item(5).RowHeight = 26.25

# Move the active selection to G4, matching the saved view state.
$ws.Range("G4").Select() | Out-Null
